$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 130960607
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 446240
$ws.Range("R5").Value = 6759818
$ws.Range("Z5").Value = "10:26"
$ws.Range("AB5").Value = "10:26"
$ws.Range("AC5").ClearContents()

# Row 6
$ws.Range("A6").Value = 130963816
$ws.Range("B6").Value = 79245
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 445932
$ws.Range("R6").Value = 6760079
$ws.Range("Z6").Value = "14:08"
$ws.Range("AB6").Value = "14:08"
$ws.Range("AC6").Value = "Rikligt i närområdet"

# Row 18
$ws.Range("A18").Value = 130960789
$ws.Range("Q18").Value = 446284
$ws.Range("R18").Value = 6759886

# Row 19
$ws.Range("A19").Value = 130960843
$ws.Range("Q19").Value = 446247
$ws.Range("R19").Value = 6759903

# Row 23
$ws.Range("A23").Value = 130963976
$ws.Range("B23").Value = 79245
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("Q23").Value = 445929
$ws.Range("R23").Value = 6760099
$ws.Range("Z23").Value = "14:08"
$ws.Range("AB23").Value = "14:08"
$ws.Range("AC23").Value = "Miljöbild"

# Row 24
$ws.Range("A24").Value = 130962640
$ws.Range("B24").Value = 79864
$ws.Range("E24").Value = 6453
$ws.Range("F24").Value = "Vedskivlav"
$ws.Range("G24").Value = "Hertelidea botryosa"
$ws.Range("H24").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q24").Value = 446038
$ws.Range("R24").Value = 6759945
$ws.Range("Z24").Value = "10:26"
$ws.Range("AB24").Value = "10:26"
$ws.Range("AC24").ClearContents()

# Row 25
$ws.Range("A25").Value = 130961746
$ws.Range("B25").Value = 57881
$ws.Range("E25").Value = 100049
$ws.Range("F25").Value = "Spillkråka"
$ws.Range("G25").Value = "Dryocopus martius"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("M25").Value = "färska spår"
$ws.Range("Q25").Value = 446098
$ws.Range("R25").Value = 6760061

# Row 26
$ws.Range("A26").Value = 130962090
$ws.Range("B26").Value = 79245
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("M26").ClearContents()
$ws.Range("Q26").Value = 446080
$ws.Range("R26").Value = 6759960

# Row 27
$ws.Range("A27").Value = 130963807
$ws.Range("B27").Value = 57881
$ws.Range("E27").Value = 100049
$ws.Range("F27").Value = "Spillkråka"
$ws.Range("G27").Value = "Dryocopus martius"
$ws.Range("H27").Value = "(Linnaeus, 1758)"
$ws.Range("M27").Value = "färska spår"
$ws.Range("Q27").Value = 445932
$ws.Range("R27").Value = 6760079
$ws.Range("Z27").Value = "14:08"
$ws.Range("AB27").Value = "14:08"

# Row 28
$ws.Range("A28").Value = 130961461
$ws.Range("Q28").Value = 446088
$ws.Range("R28").Value = 6760088
$ws.Range("AC28").ClearContents()

# Row 29
$ws.Range("A29").Value = 130961750
$ws.Range("B29").Value = 79245
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("M29").ClearContents()
$ws.Range("Q29").Value = 446098
$ws.Range("R29").Value = 6760061
$ws.Range("Z29").Value = "10:26"
$ws.Range("AB29").Value = "10:26"
$ws.Range("AC29").Value = "Rikligt i en radie av ca 50 meter"
